# Apply updated odds values to Sheet1 per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5

# Row 5
$ws.Range("G5").Value = 2.2
$ws.Range("I5").Value = 3.6
$ws.Range("L5").Value = 4.5
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("AH5").Value = 17
$ws.Range("AK5").Value = 41
$ws.Range("AO5").Value = 13

# Row 6
$ws.Range("G6").Value = 2.3
$ws.Range("I6").Value = 3.5
$ws.Range("J6").Value = 3.2
$ws.Range("S6").Value = 1.62
$ws.Range("T6").Value = 2.2
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.62
$ws.Range("AF6").Value = 81
$ws.Range("AG6").Value = 7.5
$ws.Range("AL6").Value = 51
$ws.Range("AO6").Value = 15
$ws.Range("AT6").Value = 2.2
$ws.Range("BB6").Value = 401

# Row 9
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
